$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 119, shifting existing rows 119:163 down to 120:164.
$ws.Rows.Item(119).Insert()

# Populate the newly inserted row 119 with the new record.
$ws.Cells.Item(119, 1).Value = 4
$ws.Cells.Item(119, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(119, 3).Value = "Los Lagos"
$ws.Cells.Item(119, 4).Value = 45146
$ws.Cells.Item(119, 5).Value = 10
$ws.Cells.Item(119, 6).Value = 100112031
$ws.Cells.Item(119, 7).Value = "Poroto verde"
$ws.Cells.Item(119, 8).Value = "Magnum"
$ws.Cells.Item(119, 9).Value = "Primera"
$ws.Cells.Item(119, 10).Value = 45
$ws.Cells.Item(119, 11).Value = 38000
$ws.Cells.Item(119, 12).Value = 42000
$ws.Cells.Item(119, 13).Value = 39778
$ws.Cells.Item(119, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(119, 15).Value = "Perú"
$ws.Cells.Item(119, 16).Value = 1591
$ws.Cells.Item(119, 17).Value = 25
$ws.Cells.Item(119, 18).Value = "Hortaliza"
